$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 18 and 19 swap their entire record content (Id 130983959/130983928),
# and rows 45-49 cyclically rotate their record content.
# Apply the exact per-cell before->after deltas derived from the target diff.

$ws.Range("A18").Value = 130983928
$ws.Range("AC18").Value = "Fyndplats i bitvis flerskiktad gammal granskog med frekvent förekomst av murknande björkhögstubbar för talltitans bohål."
$ws.Range("AF18").ClearContents() | Out-Null
$ws.Range("AJ18").ClearContents() | Out-Null
$ws.Range("AK18").ClearContents() | Out-Null
$ws.Range("AM18").ClearContents() | Out-Null
$ws.Range("AO18").ClearContents() | Out-Null
$ws.Range("B18").Value = 58043
$ws.Range("E18").Value = 103021
$ws.Range("F18").Value = "Talltita"
$ws.Range("G18").Value = "Poecile montanus"
$ws.Range("H18").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I18").Value = "'1"
$ws.Range("J18").ClearContents() | Out-Null
$ws.Range("L18").Value = ""
$ws.Range("M18").Value = "förbiflygande"
$ws.Range("N18").Value = "observerad"
$ws.Range("Q18").Value = 513873
$ws.Range("R18").Value = 7097763
$ws.Range("A19").Value = 130983959
$ws.Range("AC19").ClearContents() | Out-Null
$ws.Range("AF19").Value = ""
$ws.Range("AJ19").Value = "gran"
$ws.Range("AK19").Value = "Picea abies"
$ws.Range("AM19").Value = "Bark på levande träd"
$ws.Range("AO19").Value = "Bark on living woody plant # Picea abies"
$ws.Range("B19").Value = 78255
$ws.Range("E19").Value = 228579
$ws.Range("F19").Value = "Liten svartspik"
$ws.Range("G19").Value = "Chaenothecopsis nana"
$ws.Range("H19").Value = "Tibell"
$ws.Range("I19").Value = ""
$ws.Range("J19").Value = ""
$ws.Range("L19").ClearContents() | Out-Null
$ws.Range("M19").ClearContents() | Out-Null
$ws.Range("N19").Value = ""
$ws.Range("Q19").Value = 513749
$ws.Range("R19").Value = 7097845
$ws.Range("A45").Value = 130983933
$ws.Range("AJ45").Value = "sälg"
$ws.Range("AK45").Value = "Salix caprea"
$ws.Range("AM45").Value = "Bark på levande träd"
$ws.Range("AO45").Value = "Bark on living woody plant # Salix caprea"
$ws.Range("B45").Value = 80348
$ws.Range("E45").Value = 6458
$ws.Range("F45").Value = "Lunglav"
$ws.Range("G45").Value = "Lobaria pulmonaria"
$ws.Range("H45").Value = "(L.) Hoffm."
$ws.Range("K45").Value = ""
$ws.Range("Q45").Value = 514081
$ws.Range("R45").Value = 7097679
$ws.Range("A46").Value = 130982600
$ws.Range("AC46").Value = "Växer på en stående död gran."
$ws.Range("AJ46").Value = "gran"
$ws.Range("AK46").Value = "Picea abies"
$ws.Range("AM46").Value = "Stående död trädstam/högstubbe"
$ws.Range("AO46").Value = "Standing dead tree/snags # Picea abies"
$ws.Range("B46").Value = 79243
$ws.Range("E46").Value = 6425
$ws.Range("F46").Value = "Garnlav"
$ws.Range("G46").Value = "Alectoria sarmentosa"
$ws.Range("H46").Value = "(Ach.) Ach."
$ws.Range("Q46").Value = 513799
$ws.Range("R46").Value = 7098152
$ws.Range("A47").Value = 130983918
$ws.Range("AC47").ClearContents() | Out-Null
$ws.Range("AM47").Value = "Bark på levande träd"
$ws.Range("AO47").Value = "Bark on living woody plant # Picea abies"
$ws.Range("B47").Value = 83223
$ws.Range("E47").Value = 6440
$ws.Range("F47").Value = "Vitgrynig nållav"
$ws.Range("G47").Value = "Chaenotheca subroscida"
$ws.Range("H47").Value = "(Eitner) Zahlbr."
$ws.Range("Q47").Value = 513544
$ws.Range("R47").Value = 7098096
$ws.Range("A48").Value = 130983938
$ws.Range("AJ48").ClearContents() | Out-Null
$ws.Range("AK48").ClearContents() | Out-Null
$ws.Range("AM48").ClearContents() | Out-Null
$ws.Range("AO48").ClearContents() | Out-Null
$ws.Range("B48").Value = 80348
$ws.Range("E48").Value = 6458
$ws.Range("F48").Value = "Lunglav"
$ws.Range("G48").Value = "Lobaria pulmonaria"
$ws.Range("H48").Value = "(L.) Hoffm."
$ws.Range("K48").Value = "med soral"
$ws.Range("Q48").Value = 513777
$ws.Range("R48").Value = 7097981
$ws.Range("A49").Value = 130983960
$ws.Range("AJ49").Value = "gran"
$ws.Range("AK49").Value = "Picea abies"
$ws.Range("AO49").Value = "Picea abies"
$ws.Range("B49").Value = 91804
$ws.Range("E49").Value = 1108
$ws.Range("F49").Value = "Harticka"
$ws.Range("G49").Value = "Pelloporus leporinus"
$ws.Range("H49").Value = "(Fr.) Krieglst."
$ws.Range("K49").Value = "teleomorf"
$ws.Range("Q49").Value = 513513
$ws.Range("R49").Value = 7098096
